$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K (dHeight), shifting dHeight/NodeType/pattern/
# dDLBWMHz/PilotPower one column to the right (L/M/N/O/P).
$ws.Columns("K:K").Insert()

# New column K is dULCarrierMHz, holding the same values as dDLCarrierMHz (J).
$ws.Range("K1").Value = "dULCarrierMHz"
$ws.Range("K2").Value = $ws.Range("J2").Value2
$ws.Range("K3").Value = $ws.Range("J3").Value2
$ws.Range("K4").Value = $ws.Range("J4").Value2
$ws.Range("K5").Value = $ws.Range("J5").Value2

# The old "NodeType" column (now shifted to M) is renamed to "BSCat"; its
# Macro1/Pico values stay the same.
$ws.Range("M1").Value = "BSCat"

# Update the selection to match the saved workbook state.
[void]$ws.Range("B10").Select()
